$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Convert cells to "no data" text placeholders (match style/type of A14, which already uses this pattern) ---
function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range("A14").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

Set-TextCell "C15" "0"
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"
Set-TextCell "D22" "0"
Set-TextCell "E22" "***.*"
Set-TextCell "C27" "0"
Set-TextCell "D27" "0"
Set-TextCell "E27" "***.*"

# --- Numeric value updates ---
$ws.Range("N14").Value = -82.5
$ws.Range("M15").Value = -21.052631578947
$ws.Range("N15").Value = -70
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -13.333333333333
$ws.Range("I16").Value = 136
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = 4.615384615384
$ws.Range("L16").Value = -8.108108108108
$ws.Range("M16").Value = -51.601423487544
$ws.Range("N16").Value = -90.130624092888
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 18.421052631578
$ws.Range("I17").Value = 330
$ws.Range("J17").Value = 302
$ws.Range("K17").Value = 9.271523178807
$ws.Range("L17").Value = 12.244897959183
$ws.Range("M17").Value = 7.491856677524
$ws.Range("N17").Value = -58.646616541353
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -53.333333333333
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = -23.529411764705
$ws.Range("L18").Value = -37.724550898203
$ws.Range("M18").Value = -65.676567656765
$ws.Range("N18").Value = -87.850467289719
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 53.333333333333
$ws.Range("I19").Value = 293
$ws.Range("J19").Value = 267
$ws.Range("K19").Value = 9.737827715355
$ws.Range("L19").Value = -9.846153846153
$ws.Range("M19").Value = 2.090592334494
$ws.Range("N19").Value = -51.570247933884
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 15
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 83
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = 9.210526315789
$ws.Range("L20").Value = -1.190476190476
$ws.Range("M20").Value = 2.469135802469
$ws.Range("N20").Value = -77.747989276139
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 36
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 124
$ws.Range("H21").Value = 9.677419354838
$ws.Range("I21").Value = 968
$ws.Range("J21").Value = 928
$ws.Range("K21").Value = 4.310344827586
$ws.Range("L21").Value = -7.190795781399
$ws.Range("M21").Value = -24.844720496894
$ws.Range("N21").Value = -76.390243902439
$ws.Range("L22").Value = -5.882352941176
$ws.Range("C23").Value = 7
$ws.Range("E23").Value = 133.333333333333
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 76.923076923076
$ws.Range("I23").Value = 155
$ws.Range("J23").Value = 158
$ws.Range("K23").Value = -1.898734177215
$ws.Range("L23").Value = -8.284023668639
$ws.Range("M23").Value = 4.026845637583
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 138
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = 10.4
$ws.Range("I24").Value = 1015
$ws.Range("J24").Value = 959
$ws.Range("K24").Value = 5.839416058394
$ws.Range("L24").Value = -1.932367149758
$ws.Range("M24").Value = 50.593471810089
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -13.333333333333
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = -22.727272727272
$ws.Range("I25").Value = 483
$ws.Range("J25").Value = 429
$ws.Range("K25").Value = 12.587412587412
$ws.Range("L25").Value = 20.149253731343
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 49
$ws.Range("H26").Value = -18.367346938775
$ws.Range("I26").Value = 394
$ws.Range("J26").Value = 448
$ws.Range("K26").Value = -12.053571428571
$ws.Range("L26").Value = -11.261261261261
$ws.Range("M26").Value = -41.802067946824
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 13.888888888888
$ws.Range("L28").Value = 24.242424242424
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -23.809523809523
$ws.Range("M29").Value = -63.636363636363
$ws.Range("N29").Value = -91.919191919191
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 14
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -17.647058823529
$ws.Range("M30").Value = -58.823529411764
$ws.Range("N30").Value = -92
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = -83.333333333333
$ws.Range("I31").Value = 15
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = 25
$ws.Range("L31").Value = 200
$ws.Range("F33").Value = 2
$ws.Range("H33").Value = 100
$ws.Range("I33").Value = 2
$ws.Range("K33").Value = 100
$ws.Range("L33").Value = 0
